$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Upgrade OCV-SOC P25 data: update OCV value at SOC=0 from 2.7 to 2.5
$ws.Range("B2").Value = 2.5

# Refresh the embedded chart so its cached series values follow the new cell data
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$chart.Refresh()

# Update the active selection to match the new cell location
$ws.Range("S10").Select()
